# Apply crypto price/volume updates per commit "Updated cryptos list on Mon May  8 04:29:16 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.265.43"
$ws.Range("E2").Value = "  -2.29%  "

# Row 3
$ws.Range("D3").Value = "1.867.06"
$ws.Range("E3").Value = "  -1.72%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "'319.13"
$ws.Range("E5").Value = "  -1.60%  "

# Row 6
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").Value = "'0.4377"
$ws.Range("E7").Value = "  -4.54%  "

# Row 8
$ws.Range("D8").Value = "'0.3709"
$ws.Range("E8").Value = "  -2.65%  "

# Row 9
$ws.Range("D9").Value = "'0.07532"
$ws.Range("E9").Value = "  -2.26%  "

# Row 10
$ws.Range("D10").Value = "'0.9418"
$ws.Range("E10").Value = "  -3.38%  "

# Row 11
$ws.Range("D11").Value = "'21.42"
$ws.Range("E11").Value = "  -2.67%  "

# Row 12
$ws.Range("D12").Value = "1.853.15"
$ws.Range("E12").Value = "  -1.68%  "

# Row 13
$ws.Range("D13").Value = "'6.737"
$ws.Range("E13").Value = "  -2.65%  "

# Row 14
$ws.Range("D14").Value = "'5.463"
$ws.Range("E14").Value = "  -2.98%  "

# Row 15
$ws.Range("D15").Value = "'0.06865"
$ws.Range("E15").Value = "  -2.38%  "

# Row 16
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.03%  "

# Row 17
$ws.Range("D17").Value = "'82.45"
$ws.Range("E17").Value = "  -1.25%  "

# Row 18
$ws.Range("D18").Value = "'0.000009119"
$ws.Range("E18").Value = "  -3.68%  "

# Row 19
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
$ws.Range("D20").Value = "'16.01"
$ws.Range("E20").Value = "  -3.44%  "

# Row 21
$ws.Range("D21").Value = "28.240.44"
$ws.Range("E21").Value = "  -2.30%  "

# Row 22
$ws.Range("D22").Value = "'5.139"
$ws.Range("E22").Value = "  -2.70%  "

# Row 23
$ws.Range("D23").Value = "'10.84"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").Value = "2.084.73"
$ws.Range("E24").Value = "  -1.13%  "

# Row 25
$ws.Range("D25").Value = "'2.023"
$ws.Range("E25").Value = "  -3.28%  "

# Row 26
$ws.Range("D26").Value = "'154.90"
$ws.Range("E26").Value = "  -1.89%  "

# Row 27
$ws.Range("D27").Value = "'18.43"

# Row 28
$ws.Range("D28").Value = "'5.336"
$ws.Range("E28").Value = "  -4.80%  "

# Row 29
$ws.Range("D29").Value = "'114.16"
$ws.Range("E29").Value = "  -2.63%  "

# Row 30
$ws.Range("D30").Value = "'1.733"
$ws.Range("E30").Value = "  -5.29%  "

# Row 31
$ws.Range("D31").Value = "'0.09041"
$ws.Range("E31").Value = "  -2.07%  "

# Row 32
$ws.Range("D32").Value = "'0.8033"
$ws.Range("E32").Value = "  -6.37%  "

# Row 33
$ws.Range("D33").Value = "'4.858"
$ws.Range("E33").Value = "  -4.27%  "

# Row 34
$ws.Range("D34").Value = "'1.174"
$ws.Range("E34").Value = "  -5.05%  "

# Row 35
$ws.Range("D35").Value = "'2.960"
$ws.Range("E35").Value = "  -0.95%  "

# Row 37
$ws.Range("D37").Value = "'1.118"
$ws.Range("E37").Value = "  -1.66%  "

# Row 38
$ws.Range("D38").Value = "'0.05456"
$ws.Range("E38").Value = "  -3.57%  "

# Row 39
$ws.Range("D39").Value = "'0.01957"
$ws.Range("E39").Value = "  -3.36%  "

# Row 40
$ws.Range("D40").Value = "'2.966"
$ws.Range("E40").Value = "  +7.74%  "

# Row 41
$ws.Range("D41").Value = "'7.141"

# Row 42
$ws.Range("D42").Value = "'0.5256"
$ws.Range("E42").Value = "  -3.83%  "

# Row 43
$ws.Range("D43").Value = "'0.1676"
$ws.Range("E43").Value = "  -4.11%  "

# Row 44
$ws.Range("D44").Value = "'8.724"
$ws.Range("E44").Value = "  -5.64%  "

# Row 45
$ws.Range("E45").Value = "  -0.41%  "

# Row 46
$ws.Range("D46").Value = "'2.052"
$ws.Range("E46").Value = "  -0.67%  "

# Row 47
$ws.Range("D47").Value = "'0.4872"
$ws.Range("E47").Value = "  -5.21%  "

# Row 48
$ws.Range("D48").Value = "'0.000002528"
$ws.Range("E48").Value = "  -3.11%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'108.00"
$ws.Range("E49").Value = "  -1.70%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'10.53"
$ws.Range("E50").Value = "  -6.00%  "

# Row 51
$ws.Range("E51").Value = "  -4.68%  "
